# Refresh the cryptos list: updated Price / Volume(1h) figures, two rows
# where the leaderboard position (and thus the Coin/Link) changed (rows
# 20/21 and 43/44), and one coin swapped out for another (row 51,
# BabyDogeCoin -> Cronos), as pulled by the upstream GitHub Actions job.
#
# $updates holds one entry per changed row; only the columns that actually
# changed are present (B/C - Coin name/Link - are only listed for the rows
# whose coin identity changed).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row=2; D="72.625.00"; E="  +0.06%  " },
    @{ Row=3; D="2.668.28"; E="  +1.33%  " },
    @{ Row=4; E="  -0.04%  " },
    @{ Row=5; D="597.92"; E="  -1.23%  " },
    @{ Row=6; D="175.75"; E="  -1.90%  " },
    @{ Row=7; E="  -0.04%  " },
    @{ Row=8; D="0.525"; E="  -0.47%  " },
    @{ Row=9; D="2.668.43"; E="  +1.38%  " },
    @{ Row=10; E="  -3.16%  " },
    @{ Row=11; D="0.170"; E="  +2.18%  " },
    @{ Row=12; D="0.356"; E="  +0.08%  " },
    @{ Row=13; E="  -0.75%  " },
    @{ Row=14; D="3.156.18"; E="  +1.55%  " },
    @{ Row=15; D="0.0000187"; E="  -1.85%  " },
    @{ Row=16; D="72.449.43"; E="  -0.11%  " },
    @{ Row=17; D="26.29"; E="  -1.87%  " },
    @{ Row=18; D="2.649.87"; E="  +0.69%  " },
    @{ Row=19; D="12.41" },
    @{ Row=20; B="Uniswap"; C="https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"; D="8.21"; E="  +3.52%  " },
    @{ Row=21; B="BitcoinCash"; C="https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"; D="372.14"; E="  -3.50%  " },
    @{ Row=22; D="4.20"; E="  +0.16%  " },
    @{ Row=23; D="2.10"; E="  +1.74%  " },
    @{ Row=24; D="72.09"; E="  -2.65%  " },
    @{ Row=25; E="  +0.02%  " },
    @{ Row=26; D="4.35"; E="  -1.53%  " },
    @{ Row=27; D="9.88"; E="  -1.34%  " },
    @{ Row=28; D="2.789.18"; E="  +2.17%  " },
    @{ Row=29; D="1.01"; E="  +0.52%  " },
    @{ Row=30; D="0.0₃0975"; E="  +1.43%  " },
    @{ Row=31; D="8.15"; E="  +0.26%  " },
    @{ Row=32; D="497.28"; E="  -4.67%  " },
    @{ Row=33; E="  -2.05%  " },
    @{ Row=34; D="1.83"; E="  -0.18%  " },
    @{ Row=35; E="  -0.02%  " },
    @{ Row=36; D="162.14"; E="  -0.50%  " },
    @{ Row=37; D="19.56"; E="  +0.61%  " },
    @{ Row=38; D="0.113"; E="  +0.86%  " },
    @{ Row=39; D="18.94"; E="  -0.84%  " },
    @{ Row=40; D="1.39"; E="  -1.71%  " },
    @{ Row=41; D="1.77"; E="  -4.84%  " },
    @{ Row=43; B="RenderToken"; C="https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"; D="5.01"; E="  -2.93%  " },
    @{ Row=44; B="dogwifhat"; C="https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"; D="2.60"; E="  +0.77%  " },
    @{ Row=45; D="0.334"; E="  -0.20%  " },
    @{ Row=46; D="156.72"; E="  +3.83%  " },
    @{ Row=47; D="39.37"; E="  -0.19%  " },
    @{ Row=48; D="3.74"; E="  +1.05%  " },
    @{ Row=49; D="0.559"; E="  +2.32%  " },
    @{ Row=50; D="1.73"; E="  +1.93%  " },
    @{ Row=51; B="Cronos"; C="https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"; D="0.0757"; E="  -1.45%  " }
)

# Matches plain decimal numbers (e.g. "597.92", "0.0000187"). Values like
# "72.625.00" (two dots) or "0.0₃0975" (subscript digit) do NOT match, and
# are left alone since Excel already stores/treats them as text.
$numberPattern = '^-?[0-9]+(\.[0-9]+)?$'

foreach ($u in $updates) {
    $row = $u.Row

    if ($u.ContainsKey("B")) { $ws.Range("B$row").Value = $u.B }
    if ($u.ContainsKey("C")) { $ws.Range("C$row").Value = $u.C }

    if ($u.ContainsKey("D")) {
        $cell = $ws.Range("D$row")
        $text = $u.D

        if ($text -match $numberPattern) {
            # The Price column holds plain text even when the text looks like
            # a plain number (e.g. "597.92"). A bare assignment would let
            # Excel auto-convert it into a Number cell, so force text entry
            # with a leading apostrophe, then reset the resulting "quote
            # prefix" cell style back to Normal so no formatting residue is
            # left on the cell.
            $cell.Value = "'" + $text
            $cell.Style = "Normal"
        } else {
            $cell.Value = $text
        }
    }

    if ($u.ContainsKey("E")) { $ws.Range("E$row").Value = $u.E }
}
